# Actualiza horarios Linea 141 (LP1912 / LP1912-215 / 6203-6173) al corte de 10:45:47
# Refleja reordenamiento de filas por hora de scrapeo y el agregado de nuevas llegadas.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('LP1912')
$ws1.Cells.Item(2,1).Value = 'Última actualización: 10:45:47'
$ws1.Cells.Item(3,1).Value = 'Total filas: 153'
$ws1.Cells.Item(38,1).Value = '05:59:00'
$ws1.Cells.Item(38,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(38,4).Value = 66
$ws1.Cells.Item(39,1).Value = '06:50:53'
$ws1.Cells.Item(39,3).Value = '15_ABASTO'
$ws1.Cells.Item(39,4).Value = 15
$ws1.Cells.Item(53,1).Value = '06:50:53'
$ws1.Cells.Item(53,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(53,4).Value = 42
$ws1.Cells.Item(54,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(55,1).Value = '05:59:00'
$ws1.Cells.Item(55,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(55,4).Value = 93
$ws1.Cells.Item(85,1).Value = '08:22:49'
$ws1.Cells.Item(85,3).Value = '17_ROMERO'
$ws1.Cells.Item(85,4).Value = 32
$ws1.Cells.Item(86,1).Value = '08:54:22'
$ws1.Cells.Item(86,3).Value = '10_OLMOS'
$ws1.Cells.Item(86,4).Value = 0
$ws1.Cells.Item(100,1).Value = '08:54:22'
$ws1.Cells.Item(100,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(100,4).Value = 29
$ws1.Cells.Item(102,1).Value = '07:35:06'
$ws1.Cells.Item(102,3).Value = '17_ROMERO'
$ws1.Cells.Item(102,4).Value = 108
$ws1.Cells.Item(120,1).Value = '09:35:44'
$ws1.Cells.Item(120,3).Value = '17_ROMERO'
$ws1.Cells.Item(120,4).Value = 47
$ws1.Cells.Item(121,1).Value = '08:54:22'
$ws1.Cells.Item(121,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(121,4).Value = 88
$ws1.Cells.Item(129,1).Value = '10:45:47'
$ws1.Cells.Item(129,2).Value = '10:45'
$ws1.Cells.Item(129,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(129,4).Value = 0
$ws1.Cells.Item(130,1).Value = '10:45:47'
$ws1.Cells.Item(130,2).Value = '10:46'
$ws1.Cells.Item(130,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(130,4).Value = 1
$ws1.Cells.Item(131,1).Value = '10:45:47'
$ws1.Cells.Item(131,2).Value = '10:52'
$ws1.Cells.Item(131,3).Value = '15_ABASTO'
$ws1.Cells.Item(131,4).Value = 7
$ws1.Cells.Item(132,1).Value = '10:45:47'
$ws1.Cells.Item(132,2).Value = '10:53'
$ws1.Cells.Item(132,3).Value = '10_OLMOS'
$ws1.Cells.Item(132,4).Value = 8
$ws1.Cells.Item(133,1).Value = '10:45:47'
$ws1.Cells.Item(133,2).Value = '10:56'
$ws1.Cells.Item(133,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(133,4).Value = 11
$ws1.Cells.Item(134,1).Value = '10:45:47'
$ws1.Cells.Item(134,2).Value = '11:01'
$ws1.Cells.Item(134,3).Value = '215C_EL PATO'
$ws1.Cells.Item(134,4).Value = 16
$ws1.Cells.Item(135,1).Value = '10:45:47'
$ws1.Cells.Item(135,2).Value = '11:03'
$ws1.Cells.Item(135,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(135,4).Value = 18
$ws1.Cells.Item(136,1).Value = '10:45:47'
$ws1.Cells.Item(136,2).Value = '11:04'
$ws1.Cells.Item(136,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(136,4).Value = 19
$ws1.Cells.Item(136,5).Value = 'LP1912'
$ws1.Cells.Item(137,1).Value = '10:45:47'
$ws1.Cells.Item(137,2).Value = '11:06'
$ws1.Cells.Item(137,3).Value = '16_P MOR-167 Y 521'
$ws1.Cells.Item(137,4).Value = 21
$ws1.Cells.Item(137,5).Value = 'LP1912'
$ws1.Cells.Item(138,1).Value = '10:45:47'
$ws1.Cells.Item(138,2).Value = '11:11'
$ws1.Cells.Item(138,3).Value = '10_OLMOS'
$ws1.Cells.Item(138,4).Value = 26
$ws1.Cells.Item(138,5).Value = 'LP1912'
$ws1.Cells.Item(139,1).Value = '10:45:47'
$ws1.Cells.Item(139,2).Value = '11:12'
$ws1.Cells.Item(139,3).Value = '15_ABASTO'
$ws1.Cells.Item(139,4).Value = 27
$ws1.Cells.Item(139,5).Value = 'LP1912'
$ws1.Cells.Item(140,1).Value = '10:45:47'
$ws1.Cells.Item(140,2).Value = '11:19'
$ws1.Cells.Item(140,3).Value = '86_EST CHICA-ESC AGRARIA'
$ws1.Cells.Item(140,4).Value = 34
$ws1.Cells.Item(140,5).Value = 'LP1912'
$ws1.Cells.Item(141,1).Value = '10:45:47'
$ws1.Cells.Item(141,2).Value = '11:21'
$ws1.Cells.Item(141,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(141,4).Value = 36
$ws1.Cells.Item(141,5).Value = 'LP1912'
$ws1.Cells.Item(142,1).Value = '10:45:47'
$ws1.Cells.Item(142,2).Value = '11:27'
$ws1.Cells.Item(142,3).Value = '225_C ROCA-H SUR'
$ws1.Cells.Item(142,4).Value = 42
$ws1.Cells.Item(142,5).Value = 'LP1912'
$ws1.Cells.Item(143,1).Value = '10:45:47'
$ws1.Cells.Item(143,2).Value = '11:32'
$ws1.Cells.Item(143,3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(143,4).Value = 47
$ws1.Cells.Item(143,5).Value = 'LP1912'
$ws1.Cells.Item(144,1).Value = '10:45:47'
$ws1.Cells.Item(144,2).Value = '11:34'
$ws1.Cells.Item(144,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(144,4).Value = 49
$ws1.Cells.Item(144,5).Value = 'LP1912'
$ws1.Cells.Item(145,1).Value = '10:45:47'
$ws1.Cells.Item(145,2).Value = '11:35'
$ws1.Cells.Item(145,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(145,4).Value = 50
$ws1.Cells.Item(145,5).Value = 'LP1912'
$ws1.Cells.Item(146,1).Value = '10:45:47'
$ws1.Cells.Item(146,2).Value = '11:41'
$ws1.Cells.Item(146,3).Value = '17_ROMERO'
$ws1.Cells.Item(146,4).Value = 56
$ws1.Cells.Item(146,5).Value = 'LP1912'
$ws1.Cells.Item(147,1).Value = '10:45:47'
$ws1.Cells.Item(147,2).Value = '11:43'
$ws1.Cells.Item(147,3).Value = '10_OLMOS'
$ws1.Cells.Item(147,4).Value = 58
$ws1.Cells.Item(147,5).Value = 'LP1912'
$ws1.Cells.Item(148,1).Value = '10:45:47'
$ws1.Cells.Item(148,2).Value = '11:51'
$ws1.Cells.Item(148,3).Value = '215B_EL PATO'
$ws1.Cells.Item(148,4).Value = 66
$ws1.Cells.Item(148,5).Value = 'LP1912'
$ws1.Cells.Item(149,1).Value = '10:45:47'
$ws1.Cells.Item(149,2).Value = '11:59'
$ws1.Cells.Item(149,3).Value = '225_GOMEZ'
$ws1.Cells.Item(149,4).Value = 74
$ws1.Cells.Item(149,5).Value = 'LP1912'
$ws1.Cells.Item(150,1).Value = '10:45:47'
$ws1.Cells.Item(150,2).Value = '12:02'
$ws1.Cells.Item(150,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(150,4).Value = 77
$ws1.Cells.Item(150,5).Value = 'LP1912'
$ws1.Cells.Item(151,1).Value = '10:45:47'
$ws1.Cells.Item(151,2).Value = '12:06'
$ws1.Cells.Item(151,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(151,4).Value = 81
$ws1.Cells.Item(151,5).Value = 'LP1912'
$ws1.Cells.Item(152,1).Value = '10:45:47'
$ws1.Cells.Item(152,2).Value = '12:06'
$ws1.Cells.Item(152,3).Value = '14_ABASTO'
$ws1.Cells.Item(152,4).Value = 81
$ws1.Cells.Item(152,5).Value = 'LP1912'
$ws1.Cells.Item(153,1).Value = '10:45:47'
$ws1.Cells.Item(153,2).Value = '12:16'
$ws1.Cells.Item(153,3).Value = '17_ROMERO'
$ws1.Cells.Item(153,4).Value = 91
$ws1.Cells.Item(153,5).Value = 'LP1912'
$ws1.Cells.Item(154,1).Value = '10:45:47'
$ws1.Cells.Item(154,2).Value = '12:20'
$ws1.Cells.Item(154,3).Value = '215A_EL PATO'
$ws1.Cells.Item(154,4).Value = 95
$ws1.Cells.Item(154,5).Value = 'LP1912'
$ws1.Cells.Item(155,1).Value = '10:45:47'
$ws1.Cells.Item(155,2).Value = '12:20'
$ws1.Cells.Item(155,3).Value = '14_ABASTO'
$ws1.Cells.Item(155,4).Value = 95
$ws1.Cells.Item(155,5).Value = 'LP1912'
$ws1.Cells.Item(156,1).Value = '10:45:47'
$ws1.Cells.Item(156,2).Value = '12:21'
$ws1.Cells.Item(156,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(156,4).Value = 96
$ws1.Cells.Item(156,5).Value = 'LP1912'
$ws1.Cells.Item(157,1).Value = '10:45:47'
$ws1.Cells.Item(157,2).Value = '12:36'
$ws1.Cells.Item(157,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(157,4).Value = 111
$ws1.Cells.Item(157,5).Value = 'LP1912'
$ws1.Cells.Item(158,1).Value = '10:45:47'
$ws1.Cells.Item(158,2).Value = '12:38'
$ws1.Cells.Item(158,3).Value = '17_179 Y 38'
$ws1.Cells.Item(158,4).Value = 113
$ws1.Cells.Item(158,5).Value = 'LP1912'

$ws2 = $wb.Worksheets.Item('LP1912-215')
$ws2.Cells.Item(2,1).Value = 'Última actualización: 10:45:47'
$ws2.Cells.Item(3,1).Value = 'Total filas: 22'
$ws2.Cells.Item(25,1).Value = '10:45:47'
$ws2.Cells.Item(25,4).Value = 16
$ws2.Cells.Item(26,1).Value = '10:45:47'
$ws2.Cells.Item(26,2).Value = '11:51'
$ws2.Cells.Item(26,3).Value = '215B_EL PATO'
$ws2.Cells.Item(26,4).Value = 66
$ws2.Cells.Item(26,5).Value = 'LP1912'
$ws2.Cells.Item(27,1).Value = '10:45:47'
$ws2.Cells.Item(27,2).Value = '12:20'
$ws2.Cells.Item(27,3).Value = '215A_EL PATO'
$ws2.Cells.Item(27,4).Value = 95
$ws2.Cells.Item(27,5).Value = 'LP1912'

$ws3 = $wb.Worksheets.Item('6203-6173')
$ws3.Cells.Item(2,1).Value = 'Última actualización: 10:45:47'
$ws3.Cells.Item(3,1).Value = 'Total filas: 22'
$ws3.Cells.Item(25,1).Value = '10:45:47'
$ws3.Cells.Item(25,4).Value = 9
$ws3.Cells.Item(26,1).Value = '10:45:47'
$ws3.Cells.Item(26,4).Value = 28
$ws3.Cells.Item(27,1).Value = '10:45:47'
$ws3.Cells.Item(27,2).Value = '12:04'
$ws3.Cells.Item(27,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(27,4).Value = 79
$ws3.Cells.Item(27,5).Value = 'L6173'
